$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = -20.10207950654997
$ws.Cells.Item(2, 3).Value = 1.956475460748787
$ws.Cells.Item(2, 4).Value = -20.10207950654997
$ws.Cells.Item(2, 5).Value = -20.10207950654997
$ws.Cells.Item(2, 6).Value = -20.10207950654997
$ws.Cells.Item(2, 7).Value = -20.10207950654997
$ws.Cells.Item(2, 8).Value = -20.10207950654997
$ws.Cells.Item(2, 9).Value = -20.10207950654997
$ws.Cells.Item(2, 10).Value = -20.10207950654997
$ws.Cells.Item(2, 11).Value = -20.10207950654997

$ws.Cells.Item(3, 2).Value = -20.10207950654997
$ws.Cells.Item(3, 3).Value = -20.10207950654997
$ws.Cells.Item(3, 4).Value = -20.10207950654997
$ws.Cells.Item(3, 5).Value = -20.10207950654997
$ws.Cells.Item(3, 6).Value = -20.10207950654997
$ws.Cells.Item(3, 7).Value = -20.10207950654997
$ws.Cells.Item(3, 8).Value = -20.10207950654997
$ws.Cells.Item(3, 9).Value = 1.236842357814248
$ws.Cells.Item(3, 10).Value = -20.10207950654997
$ws.Cells.Item(3, 11).Value = -20.10207950654997

$ws.Cells.Item(4, 2).Value = -20.10207950654997
$ws.Cells.Item(4, 3).Value = 1.983703158569131
$ws.Cells.Item(4, 4).Value = 1.759303158574952
$ws.Cells.Item(4, 5).Value = -20.10207950654997
$ws.Cells.Item(4, 6).Value = 3.422075383645387
$ws.Cells.Item(4, 7).Value = -20.10207950654997
$ws.Cells.Item(4, 8).Value = 1.492378597792241
$ws.Cells.Item(4, 9).Value = -20.10207950654997
$ws.Cells.Item(4, 10).Value = 0.9100655509232866
$ws.Cells.Item(4, 11).Value = -20.10207950654997

$ws.Cells.Item(5, 2).Value = -20.10207950654997
$ws.Cells.Item(5, 3).Value = 1.653907042857035
$ws.Cells.Item(5, 4).Value = -20.10207950654997
$ws.Cells.Item(5, 5).Value = -20.10207950654997
$ws.Cells.Item(5, 6).Value = -20.10207950654997
$ws.Cells.Item(5, 7).Value = 2.798217778221657
$ws.Cells.Item(5, 8).Value = -20.10207950654997
$ws.Cells.Item(5, 9).Value = -20.10207950654997
$ws.Cells.Item(5, 10).Value = -20.10207950654997
$ws.Cells.Item(5, 11).Value = -20.10207950654997

$ws.Cells.Item(6, 2).Value = -20.10207950654997
$ws.Cells.Item(6, 3).Value = -20.10207950654997
$ws.Cells.Item(6, 4).Value = -20.10207950654997
$ws.Cells.Item(6, 5).Value = -20.10207950654997
$ws.Cells.Item(6, 6).Value = -20.10207950654997
$ws.Cells.Item(6, 7).Value = -20.10207950654997
$ws.Cells.Item(6, 8).Value = -20.10207950654997
$ws.Cells.Item(6, 9).Value = -20.10207950654997
$ws.Cells.Item(6, 10).Value = -20.10207950654997
$ws.Cells.Item(6, 11).Value = -20.10207950654997

$ws.Cells.Item(7, 2).Value = 2.524048596285864
$ws.Cells.Item(7, 3).Value = -20.10207950654997
$ws.Cells.Item(7, 4).Value = -20.10207950654997
$ws.Cells.Item(7, 5).Value = -20.10207950654997
$ws.Cells.Item(7, 6).Value = -20.10207950654997
$ws.Cells.Item(7, 7).Value = -20.10207950654997
$ws.Cells.Item(7, 8).Value = -20.10207950654997
$ws.Cells.Item(7, 9).Value = -20.10207950654997
$ws.Cells.Item(7, 10).Value = -20.10207950654997
$ws.Cells.Item(7, 11).Value = -20.10207950654997

$ws.Cells.Item(8, 2).Value = -20.10207950654997
$ws.Cells.Item(8, 3).Value = -20.10207950654997
$ws.Cells.Item(8, 4).Value = -20.10207950654997
$ws.Cells.Item(8, 5).Value = -20.10207950654997
$ws.Cells.Item(8, 6).Value = -20.10207950654997
$ws.Cells.Item(8, 7).Value = -20.10207950654997
$ws.Cells.Item(8, 8).Value = -20.10207950654997
$ws.Cells.Item(8, 9).Value = -20.10207950654997
$ws.Cells.Item(8, 10).Value = -20.10207950654997
$ws.Cells.Item(8, 11).Value = -20.10207950654997

$ws.Cells.Item(9, 2).Value = 3.832691984641016
$ws.Cells.Item(9, 3).Value = -20.10207950654997
$ws.Cells.Item(9, 4).Value = -20.10207950654997
$ws.Cells.Item(9, 5).Value = -20.10207950654997
$ws.Cells.Item(9, 6).Value = -20.10207950654997
$ws.Cells.Item(9, 7).Value = -20.10207950654997
$ws.Cells.Item(9, 8).Value = -20.10207950654997
$ws.Cells.Item(9, 9).Value = -20.10207950654997
$ws.Cells.Item(9, 10).Value = -20.10207950654997
$ws.Cells.Item(9, 11).Value = -20.10207950654997

$ws.Cells.Item(10, 2).Value = -20.10207950654997
$ws.Cells.Item(10, 3).Value = -20.10207950654997
$ws.Cells.Item(10, 4).Value = -20.10207950654997
$ws.Cells.Item(10, 5).Value = -20.10207950654997
$ws.Cells.Item(10, 6).Value = -20.10207950654997
$ws.Cells.Item(10, 7).Value = -20.10207950654997
$ws.Cells.Item(10, 8).Value = -20.10207950654997
$ws.Cells.Item(10, 9).Value = 1.724613663171912
$ws.Cells.Item(10, 10).Value = -20.10207950654997
$ws.Cells.Item(10, 11).Value = 2.215538738523065

$ws.Cells.Item(11, 2).Value = -20.10207950654997
$ws.Cells.Item(11, 3).Value = -20.10207950654997
$ws.Cells.Item(11, 4).Value = -20.10207950654997
$ws.Cells.Item(11, 5).Value = 4.321926877105726
$ws.Cells.Item(11, 6).Value = -20.10207950654997
$ws.Cells.Item(11, 7).Value = 2.853088810751215
$ws.Cells.Item(11, 8).Value = -20.10207950654997
$ws.Cells.Item(11, 9).Value = -20.10207950654997
$ws.Cells.Item(11, 10).Value = -20.10207950654997
$ws.Cells.Item(11, 11).Value = 1.955305780551542

$ws.Cells.Item(12, 2).Value = -20.10207950654997
$ws.Cells.Item(12, 3).Value = -20.10207950654997
$ws.Cells.Item(12, 4).Value = -20.10207950654997
$ws.Cells.Item(12, 5).Value = -20.10207950654997
$ws.Cells.Item(12, 6).Value = -20.10207950654997
$ws.Cells.Item(12, 7).Value = -20.10207950654997
$ws.Cells.Item(12, 8).Value = -20.10207950654997
$ws.Cells.Item(12, 9).Value = -20.10207950654997
$ws.Cells.Item(12, 10).Value = -20.10207950654997
$ws.Cells.Item(12, 11).Value = -20.10207950654997

$ws.Cells.Item(13, 2).Value = -20.10207950654997
$ws.Cells.Item(13, 3).Value = -20.10207950654997
$ws.Cells.Item(13, 4).Value = -20.10207950654997
$ws.Cells.Item(13, 5).Value = -20.10207950654997
$ws.Cells.Item(13, 6).Value = -20.10207950654997
$ws.Cells.Item(13, 7).Value = -20.10207950654997
$ws.Cells.Item(13, 8).Value = -20.10207950654997
$ws.Cells.Item(13, 9).Value = -20.10207950654997
$ws.Cells.Item(13, 10).Value = 1.688457120545872
$ws.Cells.Item(13, 11).Value = 1.759537904988706

$ws.Cells.Item(14, 2).Value = -20.10207950654997
$ws.Cells.Item(14, 3).Value = -20.10207950654997
$ws.Cells.Item(14, 4).Value = 1.595965920779737
$ws.Cells.Item(14, 5).Value = -20.10207950654997
$ws.Cells.Item(14, 6).Value = -20.10207950654997
$ws.Cells.Item(14, 7).Value = -20.10207950654997
$ws.Cells.Item(14, 8).Value = -20.10207950654997
$ws.Cells.Item(14, 9).Value = -20.10207950654997
$ws.Cells.Item(14, 10).Value = -20.10207950654997
$ws.Cells.Item(14, 11).Value = 1.945887424286579

$ws.Cells.Item(15, 2).Value = -20.10207950654997
$ws.Cells.Item(15, 3).Value = -20.10207950654997
$ws.Cells.Item(15, 4).Value = 1.744834420169821
$ws.Cells.Item(15, 5).Value = -20.10207950654997
$ws.Cells.Item(15, 6).Value = -20.10207950654997
$ws.Cells.Item(15, 7).Value = -20.10207950654997
$ws.Cells.Item(15, 8).Value = -20.10207950654997
$ws.Cells.Item(15, 9).Value = -20.10207950654997
$ws.Cells.Item(15, 10).Value = -20.10207950654997
$ws.Cells.Item(15, 11).Value = -20.10207950654997

$ws.Cells.Item(16, 2).Value = -20.10207950654997
$ws.Cells.Item(16, 3).Value = -20.10207950654997
$ws.Cells.Item(16, 4).Value = -20.10207950654997
$ws.Cells.Item(16, 5).Value = -20.10207950654997
$ws.Cells.Item(16, 6).Value = -20.10207950654997
$ws.Cells.Item(16, 7).Value = -20.10207950654997
$ws.Cells.Item(16, 8).Value = -20.10207950654997
$ws.Cells.Item(16, 9).Value = -20.10207950654997
$ws.Cells.Item(16, 10).Value = 1.911476598421339
$ws.Cells.Item(16, 11).Value = -20.10207950654997

$ws.Cells.Item(17, 2).Value = -20.10207950654997
$ws.Cells.Item(17, 3).Value = 2.025126745933334
$ws.Cells.Item(17, 4).Value = 1.820124861743871
$ws.Cells.Item(17, 5).Value = -20.10207950654997
$ws.Cells.Item(17, 6).Value = -20.10207950654997
$ws.Cells.Item(17, 7).Value = -20.10207950654997
$ws.Cells.Item(17, 8).Value = 2.065369573095661
$ws.Cells.Item(17, 9).Value = 2.090538698852954
$ws.Cells.Item(17, 10).Value = 2.530301835811386
$ws.Cells.Item(17, 11).Value = -20.10207950654997

$ws.Cells.Item(18, 2).Value = -20.10207950654997
$ws.Cells.Item(18, 3).Value = -20.10207950654997
$ws.Cells.Item(18, 4).Value = -20.10207950654997
$ws.Cells.Item(18, 5).Value = -20.10207950654997
$ws.Cells.Item(18, 6).Value = -20.10207950654997
$ws.Cells.Item(18, 7).Value = -20.10207950654997
$ws.Cells.Item(18, 8).Value = 1.991157825763882
$ws.Cells.Item(18, 9).Value = 2.037251066659007
$ws.Cells.Item(18, 10).Value = 2.421929594918592
$ws.Cells.Item(18, 11).Value = -20.10207950654997

$ws.Cells.Item(19, 2).Value = -20.10207950654997
$ws.Cells.Item(19, 3).Value = -20.10207950654997
$ws.Cells.Item(19, 4).Value = 1.968624859158279
$ws.Cells.Item(19, 5).Value = -20.10207950654997
$ws.Cells.Item(19, 6).Value = -20.10207950654997
$ws.Cells.Item(19, 7).Value = -20.10207950654997
$ws.Cells.Item(19, 8).Value = 1.638426371470106
$ws.Cells.Item(19, 9).Value = 1.837961967899529
$ws.Cells.Item(19, 10).Value = -20.10207950654997
$ws.Cells.Item(19, 11).Value = -20.10207950654997

$ws.Cells.Item(20, 2).Value = -20.10207950654997
$ws.Cells.Item(20, 3).Value = 1.132352622848231
$ws.Cells.Item(20, 4).Value = 1.482870611367179
$ws.Cells.Item(20, 5).Value = -20.10207950654997
$ws.Cells.Item(20, 6).Value = 3.214304372690311
$ws.Cells.Item(20, 7).Value = -20.10207950654997
$ws.Cells.Item(20, 8).Value = 1.659160456027802
$ws.Cells.Item(20, 9).Value = 1.262737514432983
$ws.Cells.Item(20, 10).Value = -20.10207950654997
$ws.Cells.Item(20, 11).Value = 2.083660091080715

$ws.Cells.Item(21, 2).Value = -20.10207950654997
$ws.Cells.Item(21, 3).Value = 1.461794974006501
$ws.Cells.Item(21, 4).Value = -20.10207950654997
$ws.Cells.Item(21, 5).Value = -20.10207950654997
$ws.Cells.Item(21, 6).Value = 2.540702230930407
$ws.Cells.Item(21, 7).Value = 1.46147848542512
$ws.Cells.Item(21, 8).Value = -20.10207950654997
$ws.Cells.Item(21, 9).Value = -20.10207950654997
$ws.Cells.Item(21, 10).Value = -20.10207950654997
$ws.Cells.Item(21, 11).Value = -20.10207950654997

